# Update crypto price/volume table with latest scrape values and bump the
# "Hora" (hour) column from 20 to 21, per the Wed Feb 8 21:15:06 UTC 2023
# GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $addr, $text) {
    # Write $text as a literal text value (matching the existing
    # text-typed cells in this sheet) without leaving a numeric/percent
    # style attached to the cell afterwards.
    $r = $sheet.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# --- Price (D) / Volume 1h (E) updates ---------------------------------
Set-CellText $ws "D2" '326.52'
Set-CellText $ws "E2" '-1.53%'
Set-CellText $ws "D3" '44.22'
Set-CellText $ws "E3" '-1.11%'
Set-CellText $ws "D4" '5.279'
Set-CellText $ws "E4" '-4.70%'
Set-CellText $ws "D5" '0.08336'
Set-CellText $ws "E5" '2.21%'
Set-CellText $ws "D6" '1.936'
Set-CellText $ws "E6" '-5.78%'
Set-CellText $ws "D7" '0.9687'
Set-CellText $ws "E7" '-1.00%'
Set-CellText $ws "D8" '2.549'
Set-CellText $ws "E8" '-3.62%'
Set-CellText $ws "D9" '0.1123'
Set-CellText $ws "E9" '0.59%'
Set-CellText $ws "D10" '0.1888'
Set-CellText $ws "E10" '-0.55%'
Set-CellText $ws "D11" '0.09689'
Set-CellText $ws "E11" '-3.49%'
Set-CellText $ws "D12" '0.04610'
Set-CellText $ws "E12" '-2.20%'
Set-CellText $ws "D13" '0.1060'
Set-CellText $ws "E13" '0.18%'
Set-CellText $ws "D14" '0.001298'
Set-CellText $ws "E14" '0.90%'
Set-CellText $ws "D15" '0.005933'
Set-CellText $ws "E15" '0.52%'
Set-CellText $ws "D16" '3.390'
Set-CellText $ws "E16" '1.23%'
Set-CellText $ws "D17" '4.400'
Set-CellText $ws "E17" '-0.70%'
Set-CellText $ws "D18" '0.3358'
Set-CellText $ws "E18" '0.20%'
Set-CellText $ws "D19" '8.561'
Set-CellText $ws "E19" '-16.46%'
Set-CellText $ws "D20" '0.1370'
Set-CellText $ws "E20" '-1.43%'
Set-CellText $ws "D21" '0.2580'
Set-CellText $ws "E21" '0.55%'
Set-CellText $ws "D22" '0.04154'
Set-CellText $ws "E22" '1.52%'
Set-CellText $ws "D23" '0.001235'
Set-CellText $ws "E23" '-4.97%'
Set-CellText $ws "D24" '0.004408'
Set-CellText $ws "E24" '-0.20%'
Set-CellText $ws "D25" '0.0001303'
Set-CellText $ws "E25" '2.00%'
Set-CellText $ws "D26" '0.0002986'
Set-CellText $ws "E26" '-20.01%'
Set-CellText $ws "D38" '0.02680'
Set-CellText $ws "E38" '-2.03%'
Set-CellText $ws "D39" '0.05557'
Set-CellText $ws "E39" '-2.96%'
Set-CellText $ws "D40" '0.007854'
Set-CellText $ws "E40" '3.49%'
Set-CellText $ws "D41" '0.1410'
Set-CellText $ws "E41" '-1.17%'
Set-CellText $ws "D42" '0.007328'
Set-CellText $ws "E42" '-2.78%'
Set-CellText $ws "D43" '0.002115'
Set-CellText $ws "E43" '8.19%'
Set-CellText $ws "D44" '0.007885'
Set-CellText $ws "E44" '-5.16%'
Set-CellText $ws "D46" '0.00006862'
Set-CellText $ws "E46" '-2.63%'
Set-CellText $ws "D47" '0.00000000752'
Set-CellText $ws "E47" '0.43%'
Set-CellText $ws "D48" '0.003493'
Set-CellText $ws "E48" '-1.19%'
Set-CellText $ws "D49" '0.003539'
Set-CellText $ws "E49" '40.67%'
Set-CellText $ws "D50" '0.00002105'
Set-CellText $ws "E50" '0.43%'
Set-CellText $ws "D51" '0.0002005'
Set-CellText $ws "E51" '0.43%'

# --- Hora (G) column: every data row moves from "20" to "21" -----------
for ($row = 2; $row -le 51; $row++) {
    Set-CellText $ws "G$row" "21"
}
